# "tolte dispense modifica DT dati"
# The header cell B1 on Foglio1 was re-cased from "lotto" to "Lotto",
# and the active selection moved from the whole table (A1:C33) to B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")
$ws.Activate()

$ws.Range("B1").Value = "Lotto"

$ws.Range("B2").Select()
